# Group Reflection - Charles.docx
# Commit: "Group reflection update + CSS update"
#
# 1) Add a new bullet ("Originally didn't like having to use GitHub, this
#    later changed as I realised its...") under "What was surprising?",
#    directly above "Team members personalities worked well together" -
#    same ListParagraph / numId=3 bullet list.
# 2) Under "What have I learned about groups?", the trailing empty bullet
#    (ListParagraph / numId=3) that followed "Delegation of roles is
#    important from the beginning" loses its list formatting and an extra
#    blank paragraph is added, leaving five plain empty paragraphs at the
#    end of that section.

$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# Locate the "Team members personalities..." bullet and insert a new,
# identically-formatted bullet immediately before it.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Team members personalities worked well together", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $teamPara = $rng1.Paragraphs(1)
    $teamIndex = $teamPara.Index

    # Inherits the ListParagraph / numId=3 paragraph formatting automatically.
    $teamPara.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs($teamIndex)
    $newPara.Range.Text = "Originally didn’t like having to use GitHub, this later changed as I realised its…"
}

# --- Change 2 --------------------------------------------------------
# Locate "Delegation of roles is important from the beginning"; the
# paragraph right after it is the empty ListParagraph/numId=3 bullet that
# needs to become plain, with one extra plain paragraph inserted too.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Delegation of roles is important from the beginning", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $delPara = $rng2.Paragraphs(1)
    $emptyIndex = $delPara.Index + 1
    $emptyListPara = $d.Paragraphs($emptyIndex)

    # A bare <w:p/> (no pStyle/numPr) - replaces the empty bullet in place.
    $plainXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
    $null = $emptyListPara.Range.InsertXML($plainXml)

    # Add a second, brand-new plain empty paragraph right after it.
    $firstPlain = $d.Paragraphs($emptyIndex)
    $firstPlain.Range.InsertParagraphAfter()

    $secondPlain = $d.Paragraphs($emptyIndex + 1)
    $null = $secondPlain.Range.InsertXML($plainXml)
}

Write-Host "Group reflection update applied."
